# Updated cryptos list on Fri Jul 12 06:49:48 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "57.311.20"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  -1.16%  "
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "3.095.58"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  +0.04%  "
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "522.98"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "136.72"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  +0.03%  "
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "3.095.93"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  +3.33%  "
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "7.37"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("E11").Value = "  -1.50%  "
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "0.401"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  +2.26%  "
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "3.629.51"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  +0.90%  "
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "25.41"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("E16").Value = "  -1.92%  "
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "57.416.93"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  -1.04%  "
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "3.098.40"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E19").Value = "  -3.29%  "
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "12.50"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("E21").Value = "  -1.32%  "
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "350.13"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("E23").Value = "  -0.08%  "
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "68.78"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("E26").Value = "  -1.52%  "
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "0.998"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0868"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  -6.25%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -0.08%  "
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "5.86"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -7.80%  "
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "20.93"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  -0.28%  "
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "4.88"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  +5.17%  "
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "1.14"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  -3.81%  "
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "159.08"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "6.04"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -1.97%  "
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "25.64"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("E39").Value = "  -0.70%  "
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0659"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -1.74%  "
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "1.58"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("E42").Value = "  +1.09%  "
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "0.696"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  +1.34%  "
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "2.402.12"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  +5.87%  "
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "36.71"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("E46").Value = "  +0.04%  "
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "3.134.95"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("E49").Value = "  -3.66%  "
$ws.Range("E50").Value = "  -2.54%  "
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "0.766"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +1.59%  "
